# KHL stats refresh: new slate of games for 2025-11-13, updated rolling
# shots-on-goal aggregates, and a bumped Meta_ext build/version stamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Matches_SOG: append the two newly-played/scheduled games (uids 897767,
# 897768) right after the existing last row (469).
# ---------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 470
Set-TextCell $wsMatches 470 1 "897767"
$wsMatches.Cells.Item(470, 2).Value = "2025-11-13T16:30:00"
$wsMatches.Cells.Item(470, 3).Value = "Авангард"
$wsMatches.Cells.Item(470, 4).Value = "Трактор"
$wsMatches.Cells.Item(470, 5).Value = 21
$wsMatches.Cells.Item(470, 6).Value = 35
$wsMatches.Cells.Item(470, 7).Value = "khl_text"

# Row 471
Set-TextCell $wsMatches 471 1 "897768"
$wsMatches.Cells.Item(471, 2).Value = "2025-11-13T19:00:00"
$wsMatches.Cells.Item(471, 3).Value = "Ак Барс"
$wsMatches.Cells.Item(471, 4).Value = "Динамо М"
$wsMatches.Cells.Item(471, 5).Value = 33
$wsMatches.Cells.Item(471, 6).Value = 27
$wsMatches.Cells.Item(471, 7).Value = "khl_text"

# ---------------------------------------------------------------------
# Shots_HA: as_of_utc moves to the new snapshot time for every team, and
# the teams involved in today's games get their rolling HA shot tallies
# refreshed.
# ---------------------------------------------------------------------
$wsHA = $wb.Worksheets.Item("Shots_HA")

$wsHA.Range("D2:D23").Value = "2025-11-13T19:00:00Z"

# Авангард (home vs Трактор)
$wsHA.Range("E2").Value = 23
$wsHA.Range("G2").Value = 751
$wsHA.Range("H2").Value = 651
$wsHA.Range("I2").Value = 32.7
$wsHA.Range("J2").Value = 28.3

# Ак Барс (home vs Динамо М)
$wsHA.Range("E5").Value = 24
$wsHA.Range("G5").Value = 812
$wsHA.Range("H5").Value = 614
$wsHA.Range("I5").Value = 33.8
$wsHA.Range("J5").Value = 25.6

# Динамо М (away at Ак Барс)
$wsHA.Range("F8").Value = 22
$wsHA.Range("K8").Value = 623
$wsHA.Range("L8").Value = 706
$wsHA.Range("M8").Value = 28.3
$wsHA.Range("N8").Value = 32.1

# Трактор (away at Авангард)
$wsHA.Range("F21").Value = 27
$wsHA.Range("K21").Value = 923
$wsHA.Range("N21").Value = 32.1

# ---------------------------------------------------------------------
# Shots_Summary: same as_of_utc bump, plus totals for the four teams
# playing today.
# ---------------------------------------------------------------------
$wsSS = $wb.Worksheets.Item("Shots_Summary")

$wsSS.Range("D2:D23").Value = "2025-11-13T19:00:00Z"

# Авангард
$wsSS.Range("E2").Value = 41
$wsSS.Range("F2").Value = 1393
$wsSS.Range("G2").Value = 1201
$wsSS.Range("H2").Value = 34
$wsSS.Range("I2").Value = 29.3

# Ак Барс
$wsSS.Range("E5").Value = 45
$wsSS.Range("F5").Value = 1523
$wsSS.Range("G5").Value = 1241
$wsSS.Range("H5").Value = 33.8

# Динамо М
$wsSS.Range("E8").Value = 40
$wsSS.Range("F8").Value = 1221
$wsSS.Range("G8").Value = 1197
$wsSS.Range("H8").Value = 30.5
$wsSS.Range("I8").Value = 29.9

# Трактор
$wsSS.Range("E21").Value = 45
$wsSS.Range("F21").Value = 1522
$wsSS.Range("G21").Value = 1409
$wsSS.Range("I21").Value = 31.3

# ---------------------------------------------------------------------
# Meta_ext: bump the snapshot timestamp and the build_version counter.
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Range("B2").Value = "2025-11-13T19:00:00Z"
$wsMeta.Range("D2").Value = 64

Write-Output "KHL refresh applied"
